# Applies the "first set of weapons tested" edit:
#  - Adds a new "Critical" entry to Sheet1's Type list (column G), row 15
#  - Updates Sheet2's test row (row 3) selections to: All / Crit Threshold / Critical
#  - Updates sheet view selections/active tab to match the saved state

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet2 = $wb.Worksheets.Item("Sheet2")

# Add new shared value "Critical" to Sheet1 list (column G, row 15)
$sheet1.Range("G15").Value = "Critical"
$sheet1.Range("G15").Style = $sheet1.Range("G14").Style

# Update Sheet2's lookup test row with the new selections
$sheet2.Range("D3").Value = "All"
$sheet2.Range("E3").Value = "Crit Threshold"
$sheet2.Range("F3").Value = "Critical"

# Update active selections / active sheet to match saved workbook state
$sheet1.Range("F29").Select()
$sheet2.Range("F3").Select()

$sheet1.Activate()
